$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 251, shifting existing rows 251:267 down to 252:268.
$ws.Rows("251:251").Insert()

# Populate the newly inserted row 251 with the new data record.
$ws.Cells.Item(251, 1).Value = 10
$ws.Cells.Item(251, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(251, 3).Value = "La Araucanía"
$ws.Cells.Item(251, 4).Value = 44610
$ws.Cells.Item(251, 5).Value = 9
$ws.Cells.Item(251, 6).Value = 100114013
$ws.Cells.Item(251, 7).Value = "Zanahoria"
$ws.Cells.Item(251, 8).Value = "Sin especificar"
$ws.Cells.Item(251, 9).Value = "Primera"
$ws.Cells.Item(251, 10).Value = 50
$ws.Cells.Item(251, 11).Value = 8000
$ws.Cells.Item(251, 12).Value = 8000
$ws.Cells.Item(251, 13).Value = 8000
$ws.Cells.Item(251, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(251, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(251, 16).Value = 320
$ws.Cells.Item(251, 17).Value = 25
$ws.Cells.Item(251, 18).Value = "Hortaliza"

# Preserve the date-cell number format used throughout column D.
$ws.Cells.Item(251, 4).NumberFormat = $ws.Cells.Item(252, 4).NumberFormat
